$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Header block (rows 4-6): "Created:" and "Last modified:" swap
#    places, and "Last modified:" picks up a new date/description.
#    Old: row4=Created/7-16-20, row5=Last modified(blank), row6=blank
#    New: row4=Last modified/new text, row5=blank, row6=Created/7-16-20
# ---------------------------------------------------------------
$ws.Range("B4").Value = "Last modified:"
$ws.Range("D4").Value = "06/01/22 - SAC - updated 2022 look-up values w/ latest 2022 kW multipliers"
$ws.Range("D4").NumberFormat = "m/d/yyyy"
$ws.Range("E4").NumberFormat = "m/d/yyyy"

$ws.Range("B6").Value = "Created:"
$ws.Range("D6").Value = "7/16/20 - SAC"
$ws.Range("D6").NumberFormat = "General"

# ---------------------------------------------------------------
# 2) "Mod history:" row (row 9) gains a new entry describing the
#    2022 SMUD Community Solar update.
# ---------------------------------------------------------------
$ws.Range("D9").Value = "05/25/22 - SAC - added newly calculated SMUD Community Solar results for use in 2022 analysis"
$ws.Range("D9").NumberFormat = "m/d/yyyy"

# ---------------------------------------------------------------
# 3) Insert a new dependent-variable row (new row 20) describing
#    "kBtu of Source Energy (EDR1)" right after the CO2 row.
# ---------------------------------------------------------------
$ws.Rows(20).Insert()
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = "kBtu of Source Energy (EDR1)"

# ---------------------------------------------------------------
# 4) Insert a new table data row (new row 26) for the 2022 SMUD
#    Community Solar results, directly after the 2019 data row.
# ---------------------------------------------------------------
$ws.Rows(26).Insert()

# New column header for the table: "SrcKBtu"
$ws.Range("I24").Value = "SrcKBtu"
$ws.Range("I24").Copy()
$ws.Range("I24").PasteSpecial(-4122)

# 2019 data row gains a SrcKBtu value of 0
$ws.Range("I25").Value = 0

# New 2022 data row values
$ws.Range("C26").Value = 2022
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = -2137.9
$ws.Range("F26").Value = -0.0123
$ws.Range("G26").Value = -47454
$ws.Range("H26").Value = -99.151
$ws.Range("I26").Value = -1609.6
$ws.Range("J26").Value = ";"
$ws.Range("K26").Value = "SMUD Neighborhood SolarShares - Wildflower"

# Formatting for the new row: copy number formatting from the
# analogous 2019 row, then adjust font/border to match the sheet.
$ws.Range("C25:D25").Copy()
$ws.Range("C26:D26").PasteSpecial(-4122)
$ws.Range("J25").Copy()
$ws.Range("J26").PasteSpecial(-4122)

$ws.Range("D15").Copy()
$ws.Range("E26:I26").PasteSpecial(-4122)

$ws.Range("E26").Borders.Item(7).LineStyle = 1
$ws.Range("E26").Borders.Item(7).Weight = 2
$ws.Range("E26").Borders.Item(7).Color = 0

# ENDTABLE filler row ("*" row, now row 27) gains a SrcKBtu value of 0
$ws.Range("I27").Value = 0

Write-Output "edit complete"
